$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 5241.343
$ws.Range("I15").Value = 5241.343
$ws.Range("K15").Value = 15724.029
$ws.Range("M15").Value = -15555.029
# Row 22
$ws.Range("H22").Value = 1000
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3344
# Row 74
$ws.Range("H74").Value = 9203.478999999999
$ws.Range("I74").Value = 7083.636
$ws.Range("K74").Value = 7083.636
$ws.Range("M74").Value = -6147.636
# Row 77
$ws.Range("H77").Value = 9203.478999999999
$ws.Range("I77").Value = 7083.636
$ws.Range("K77").Value = 35418.18
$ws.Range("M77").Value = -30738.18
# Row 86
$ws.Range("H86").Value = 2250.6428
$ws.Range("J86").Value = 1877.4286
$ws.Range("L86").Value = 1877.4286
$ws.Range("N86").Value = -4123.4286
# Row 89
$ws.Range("H89").Value = 2250.6428
$ws.Range("J89").Value = 1877.4286
$ws.Range("L89").Value = 9387.143
$ws.Range("N89").Value = -20619.143
# Row 132
$ws.Range("H132").Value = 2076.08
$ws.Range("I132").Value = 1945.5454
$ws.Range("K132").Value = 5836.6362
$ws.Range("M132").Value = -3306.6362
# Row 135
$ws.Range("H135").Value = 683.3
$ws.Range("I135").Value = 704.7895
$ws.Range("K135").Value = 6343.1055
$ws.Range("M135").Value = -3808.1055
# Row 137
$ws.Range("H137").Value = 1146.1177
$ws.Range("I137").Value = 1170.6666
$ws.Range("K137").Value = 3511.9998
$ws.Range("M137").Value = -961.9998000000001
# Row 138
$ws.Range("H138").Value = 2958
$ws.Range("I138").Value = 2840.2856
$ws.Range("K138").Value = 8520.856800000001
$ws.Range("M138").Value = -3380.856800000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8769.405000000001
$ws.Range("I32").Value = 7135
$ws.Range("K32").Value = 7135
$ws.Range("M32").Value = -6848
# Row 37
$ws.Range("H37").Value = 17000
$ws.Range("J37").Value = 17000
$ws.Range("L37").Value = 17000
$ws.Range("N37").Value = -17546
# Row 123
$ws.Range("H123").Value = 59996.332
$ws.Range("J123").Value = 59996.332
$ws.Range("L123").Value = 59996.332
$ws.Range("N123").Value = -69796.33199999999
# Row 132
$ws.Range("H132").Value = 1313.9362
$ws.Range("I132").Value = 1318
$ws.Range("K132").Value = 3954
$ws.Range("M132").Value = -1424

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3432.5
$ws.Range("J20").Value = 3523.375
$ws.Range("L20").Value = 3523.375
$ws.Range("N20").Value = -4017.375
# Row 86
$ws.Range("H86").Value = 1206.3684
$ws.Range("I86").Value = 1135.909
$ws.Range("K86").Value = 1135.909
$ws.Range("M86").Value = -12.90900000000011
# Row 89
$ws.Range("H89").Value = 1206.3684
$ws.Range("I89").Value = 1135.909
$ws.Range("K89").Value = 5679.545
$ws.Range("M89").Value = -63.54500000000007
# Row 105
$ws.Range("H105").Value = 4687.5
$ws.Range("I105").Value = 4687.5
$ws.Range("K105").Value = 4687.5
$ws.Range("M105").Value = -2940.5
# Row 133
$ws.Range("H133").Value = 26043.479
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 134
$ws.Range("H134").Value = 2736.8857
$ws.Range("I134").Value = 2313.1453
$ws.Range("K134").Value = 6939.4359
$ws.Range("M134").Value = -4404.4359

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3762.375
$ws.Range("I16").Value = 2920
$ws.Range("J16").Value = 5166.3335
$ws.Range("K16").Value = 2920
$ws.Range("L16").Value = 5166.3335
$ws.Range("M16").Value = -2633
$ws.Range("N16").Value = -5740.3335
# Row 19
$ws.Range("H19").Value = 330.42856
$ws.Range("J19").Value = 25
$ws.Range("L19").Value = 25
$ws.Range("N19").Value = -365
# Row 24
$ws.Range("H24").Value = 330.42856
$ws.Range("J24").Value = 25
$ws.Range("L24").Value = 25
$ws.Range("N24").Value = -365
# Row 31
$ws.Range("H31").Value = 4669.3105
$ws.Range("I31").Value = 3069.5293
$ws.Range("K31").Value = 3069.5293
$ws.Range("M31").Value = -2774.5293
# Row 34
$ws.Range("H34").Value = 4669.3105
$ws.Range("I34").Value = 3069.5293
$ws.Range("K34").Value = 3069.5293
$ws.Range("M34").Value = -2867.5293
# Row 99
$ws.Range("H99").Value = 8750
$ws.Range("I99").Value = 10000
$ws.Range("K99").Value = 10000
$ws.Range("M99").Value = -8502
# Row 113
$ws.Range("H113").Value = 3762.375
$ws.Range("I113").Value = 2920
$ws.Range("J113").Value = 5166.3335
$ws.Range("K113").Value = 2920
$ws.Range("L113").Value = 5166.3335
$ws.Range("M113").Value = -750
$ws.Range("N113").Value = -9506.333500000001
# Row 126
$ws.Range("H126").Value = 8750
$ws.Range("I126").Value = 10000
$ws.Range("K126").Value = 30000
$ws.Range("M126").Value = -27530
# Row 141
$ws.Range("H141").Value = 215717.19
$ws.Range("J141").Value = 215717.19
$ws.Range("L141").Value = 215717.19
$ws.Range("N141").Value = -226077.19

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 504.0625
$ws.Range("J12").Value = 794.8
$ws.Range("L12").Value = 2384.4
$ws.Range("N12").Value = -2730.4
# Row 14
$ws.Range("H14").Value = 174.625
$ws.Range("I14").Value = 174.625
$ws.Range("K14").Value = 523.875
$ws.Range("M14").Value = -350.875
# Row 33
$ws.Range("H33").Value = 297.1875
$ws.Range("I33").Value = 144.66667
$ws.Range("J33").Value = 388.7
$ws.Range("K33").Value = 868.0000200000001
$ws.Range("L33").Value = 2332.2
$ws.Range("M33").Value = -585.0000200000001
$ws.Range("N33").Value = -2898.2
# Row 38
$ws.Range("H38").Value = 1420.5834
$ws.Range("I38").Value = 1299.6666
$ws.Range("J38").Value = 1460.8889
$ws.Range("K38").Value = 3898.9998
$ws.Range("L38").Value = 4382.6667
$ws.Range("M38").Value = -3551.9998
$ws.Range("N38").Value = -5076.6667
# Row 68
$ws.Range("H68").Value = 41670336
$ws.Range("I68").Value = 55559116
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 166677348
$ws.Range("L68").Value = 12000
$ws.Range("M68").Value = -166676537
$ws.Range("N68").Value = -13622
# Row 71
$ws.Range("H71").Value = 41670336
$ws.Range("I71").Value = 55559116
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 500032044
$ws.Range("L71").Value = 36000
$ws.Range("M71").Value = -500027988
$ws.Range("N71").Value = -44112
# Row 104
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("N104").ClearContents()
# Row 129
$ws.Range("H129").Value = 15657848
$ws.Range("I129").Value = 37148820
$ws.Range("J129").Value = 779481.9
$ws.Range("K129").Value = 111446460
$ws.Range("L129").Value = 2338445.7
$ws.Range("M129").Value = -111441460
$ws.Range("N129").Value = -2348445.7
# Row 131
$ws.Range("H131").Value = 19233304
$ws.Range("I131").Value = 71429560
$ws.Range("J131").Value = 3102.842
$ws.Range("K131").Value = 214288680
$ws.Range("L131").Value = 9308.526
$ws.Range("M131").Value = -214283640
$ws.Range("N131").Value = -19388.526

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 2849.625
$ws.Range("I14").Value = 2274.5
$ws.Range("K14").Value = 2274.5
$ws.Range("M14").Value = -2106.5
# Row 122
$ws.Range("H122").Value = 1656.7273
$ws.Range("I122").Value = 1740.5625
$ws.Range("K122").Value = 5221.6875
$ws.Range("M122").Value = -2771.6875
# Row 126
$ws.Range("H126").Value = 3499.7778
$ws.Range("I126").Value = 2636.0908
$ws.Range("K126").Value = 7908.2724
$ws.Range("M126").Value = -5438.2724
# Row 132
$ws.Range("H132").Value = 3455.8572
$ws.Range("I132").Value = 3198.2083
$ws.Range("J132").Value = 5001.75
$ws.Range("K132").Value = 9594.624899999999
$ws.Range("L132").Value = 15005.25
$ws.Range("M132").Value = -7064.624899999999
$ws.Range("N132").Value = -20065.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
# Row 40
$ws.Range("H40").Value = 4284.4443
$ws.Range("I40").Value = 4512
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 4512
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -4376
$ws.Range("N40").Value = -4272
# Row 93
$ws.Range("H93").Value = 12903.944
$ws.Range("I93").Value = 1327.3
$ws.Range("K93").Value = 1327.3
$ws.Range("M93").Value = -79.29999999999995
# Row 111
$ws.Range("H111").Value = 10000
$ws.Range("J111").Value = 10000
$ws.Range("L111").Value = 10000
$ws.Range("N111").Value = -18180
# Row 132
$ws.Range("H132").Value = 4655.2163
$ws.Range("I132").Value = 5036.7095
$ws.Range("J132").Value = 2684.1667
$ws.Range("K132").Value = 15110.1285
$ws.Range("L132").Value = 8052.500100000001
$ws.Range("M132").Value = -12580.1285
$ws.Range("N132").Value = -13112.5001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 3999
$ws.Range("J96").Value = 3999
$ws.Range("L96").Value = 3999
$ws.Range("N96").Value = -6745
# Row 100
$ws.Range("H100").Value = 561.65
$ws.Range("I100").Value = 535
$ws.Range("K100").Value = 1070
$ws.Range("M100").Value = -529
# Row 113
$ws.Range("H113").Value = 541.25
$ws.Range("I113").Value = 332.46155
$ws.Range("J113").Value = 788
$ws.Range("K113").Value = 997.38465
$ws.Range("L113").Value = 2364
$ws.Range("M113").Value = 1172.61535
$ws.Range("N113").Value = -6704
# Row 122
$ws.Range("H122").Value = 3252.8462
$ws.Range("I122").Value = 2135.4736
$ws.Range("K122").Value = 6406.4208
$ws.Range("M122").Value = -3956.4208
